# Penalty Reward System (unfinished) - data refresh for forecast_summary_B08R5736B3.xlsx
#
# Shifts the "Forecast Comparison" week-start dates forward by one week,
# writes the new MyForecast values, and refreshes the derived figures on
# the "Summary" sheet to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Forecast Comparison"
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

# New Week_Start_Date values (column B, rows 2-17) - each week shifted
# forward by one (the old W16 date 2025-04-20 now belongs to W15, and a
# new trailing date 2025-04-27 appears for W16).
$newWeekStart = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

# New MyForecast values (column D, rows 2-17).
$newForecast = @(84, 99, 103, 102, 109, 123, 86, 85, 81, 82, 77, 109, 95, 76, 61, 54)

# Keep column B as plain text (it holds text dates like "2025-01-12",
# not real Excel date serials) so the write doesn't get auto-converted
# into a date value.
$wsForecast.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 2).Value = $newWeekStart[$i]
    $wsForecast.Cells.Item($row, 4).Value = $newForecast[$i]
}

# ---------------------------------------------------------------------
# Sheet 2: "Summary"
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

# All Value cells on the Summary sheet are stored as text, even the
# numeric-looking ones - keep them text on write too.
$wsSummary.Range("B2").NumberFormat = "@"
$wsSummary.Range("B4:B6").NumberFormat = "@"
$wsSummary.Range("B8:B15").NumberFormat = "@"

$wsSummary.Range("B2").Value = "2022-12-25 to 2025-01-05"

$wsSummary.Range("B4").Value = "393"
$wsSummary.Range("B5").Value = "151"
$wsSummary.Range("B6").Value = "128"

$wsSummary.Range("B8").Value = "16036 units"
$wsSummary.Range("B9").Value = "1425"
$wsSummary.Range("B10").Value = "791"
$wsSummary.Range("B11").Value = "388"
$wsSummary.Range("B12").Value = "123"
$wsSummary.Range("B13").Value = "2025-02-16"
$wsSummary.Range("B14").Value = "54"
$wsSummary.Range("B15").Value = "2025-04-27"

Write-Output "Penalty Reward System data refresh applied."
